$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 9 ("test" item) -- add this shared string first so it gets the
# lower shared-string table index (matches authoring order).
$ws.Range("A9").Value = "test"
$ws.Range("B9").Value = 123

# New header in C1
$ws.Range("C1").Value = "分機號碼"

# New "extension number" numeric data in column C, rows 2-9 (22222 everywhere)
$ws.Range("C2").Value = 22222
$ws.Range("C3").Value = 22222
$ws.Range("C4").Value = 22222
$ws.Range("C5").Value = 22222
$ws.Range("C6").Value = 22222
$ws.Range("C7").Value = 22222
$ws.Range("C8").Value = 22222
$ws.Range("C9").Value = 22222

# New numeric data in column B where previously empty
$ws.Range("B2").Value = 123123
$ws.Range("B3").Value = 123123
$ws.Range("B5").Value = 123123
$ws.Range("B6").Value = 123123
$ws.Range("B8").Value = 123123

# Column widths to match authored layout (26.42578125 / 35.5703125 chars,
# expressed in the COM ColumnWidth scale so the stored OOXML width lands as
# close as the host's pixel-rounded width model allows)
$ws.Columns.Item(2).ColumnWidth = 25.714285714285715
$ws.Columns.Item(3).ColumnWidth = 34.857142857142854

# Restore the active selection cell used by the author
$ws.Range("C5").Select()
